# Edit slide 5 ("Jednostki w css") of the presentation:
#  - reposition the existing "vw" and "%" boxes (rounded-rect + textbox pairs)
#  - add two new boxes (rounded-rect + textbox pairs) for "em" and "rem" css units
#
# Shapes are matched by their (stable) Name, not by positional index, so the
# script is resilient to any incidental shape-order differences.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

# ---- the four existing "rounded rectangle" style shapes (no text themselves) ----
$rectVw  = Get-ShapeByName $s "Prostokąt: zaokrąglone rogi 10"   # vw
$rectPct = Get-ShapeByName $s "Prostokąt: zaokrąglone rogi 8"    # %

# ---- the four existing textbox shapes (carry the actual text) ----
$tbVw  = Get-ShapeByName $s "pole tekstowe 5"   # vw
$tbPct = Get-ShapeByName $s "pole tekstowe 6"   # %

# 1) Move the "vw" rounded-rect + textbox to their new spot (bottom-middle).
$rectVw.Left   = 342.38536133070863
$rectVw.Top    = 404.9874115748031
$rectVw.Width  = 193.5412678425197
$rectVw.Height = 72.70315260629921

$tbVw.Left   = 342.38536133070863
$tbVw.Top    = 404.87056018110235
$tbVw.Width  = 193.5412678425197
$tbVw.Height = 72.70315260629921

# 2) Move the "%" rounded-rect + textbox to their new spot (bottom-left).
$rectPct.Left   = 112.2935453070866
$rectPct.Top    = 435.2473915748032
$rectPct.Width  = 157.21109836220472
$rectPct.Height = 50.89220672440945

$tbPct.Left   = 112.2936220472441
$tbPct.Top    = 432.5966959133858
$tbPct.Width  = 157.21102162204724
$tbPct.Height = 50.89220672440945

# 3) Duplicate rounded-rect shapes for the new "em" and "rem" boxes, reusing
#    the existing shape style (fill/line/effect theme refs).
$emRectDup = $rectPct.Duplicate()
$rectEm = $emRectDup.Item(1)
$rectEm.Name   = "Prostokąt: zaokrąglone rogi 12"
$rectEm.Left   = 686.9724429448819
$rectEm.Top    = 282.38536133070863
$rectEm.Width  = 225.9082727165354
$rectEm.Height = 72.70315260629921

$remRectDup = $rectVw.Duplicate()
$rectRem = $remRectDup.Item(1)
$rectRem.Name   = "Prostokąt: zaokrąglone rogi 14"
$rectRem.Left   = 635.0090641181102
$rectRem.Top    = 432.5966041732284
$rectRem.Width  = 212.69731983464567
$rectRem.Height = 72.7032283464567

# 4) Duplicate textbox shapes for the new "em" and "rem" labels, reusing the
#    existing textbox formatting (no fill, centered, auto-fit), then set text.
$emTbDup = $tbPct.Duplicate()
$tbEm = $emTbDup.Item(1)
$tbEm.Name = "pole tekstowe 11"
$tbEm.TextFrame.TextRange.Text = "em- działa jak x, może być zmienny dla każego elementu"
$tbEm.Left   = 686.9724429448819
$tbEm.Top    = 282.38536133070863
$tbEm.Width  = 225.9082727165354
$tbEm.Height = 72.70315260629921

$remTbDup = $tbVw.Duplicate()
$tbRem = $remTbDup.Item(1)
$tbRem.Name = "pole tekstowe 13"
$tbRem.TextFrame.TextRange.Text = "rem- działa jak x, jest taki sam dla całego dokumentu"
$tbRem.Left   = 635.0090641181102
$tbRem.Top    = 432.5966959133858
$tbRem.Width  = 212.69731983464567
$tbRem.Height = 72.70315260629921
